# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (G) values computed from the regenerated save_data.
$kValues = @{
    2  = 2
    3  = 7
    4  = 8
    5  = 3
    6  = 9
    7  = 5
    8  = 6
    9  = 4
    10 = 1
    11 = 7
    12 = 3
    13 = 8
    14 = 4
    15 = 4
    16 = 3
    17 = 4
    18 = 6
    19 = 7
    20 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
